# Auto-generated Excel COM-interop edit script
# Applies value updates to the Leve profit-tracking tables across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as captured by the source diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6375.5
$ws.Range("I18").Value = 6375.5
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 6375.5
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -6091.5
$ws.Range("N18").ClearContents()
$ws.Range("H74").Value = 3443.2856
$ws.Range("I74").Value = 3999
$ws.Range("J74").Value = 3350.6667
$ws.Range("K74").Value = 3999
$ws.Range("L74").Value = 3350.6667
$ws.Range("M74").Value = -3063
$ws.Range("N74").Value = -5222.6667
$ws.Range("H77").Value = 3443.2856
$ws.Range("I77").Value = 3999
$ws.Range("J77").Value = 3350.6667
$ws.Range("K77").Value = 19995
$ws.Range("L77").Value = 16753.3335
$ws.Range("M77").Value = -15315
$ws.Range("N77").Value = -26113.3335
$ws.Range("H97").Value = 17466.75
$ws.Range("J97").Value = 17466.75
$ws.Range("L97").Value = 52400.25
$ws.Range("N97").Value = -53392.25
$ws.Range("H107").Value = 498.16666
$ws.Range("I107").Value = 476.34784
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 476.34784
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1443.65216
$ws.Range("N107").Value = -4840
$ws.Range("H116").Value = 6110.5713
$ws.Range("I116").Value = 4806.304
$ws.Range("J116").Value = 12110.2
$ws.Range("K116").Value = 4806.304
$ws.Range("L116").Value = 12110.2
$ws.Range("M116").Value = -1364.304
$ws.Range("N116").Value = -18994.2
$ws.Range("H121").Value = 1100.3334
$ws.Range("I121").Value = 301
$ws.Range("J121").Value = 1500
$ws.Range("K121").Value = 903
$ws.Range("L121").Value = 4500
$ws.Range("M121").Value = 844
$ws.Range("N121").Value = -7994
$ws.Range("H126").Value = 20000
$ws.Range("J126").Value = 20000
$ws.Range("L126").Value = 20000
$ws.Range("N126").Value = -29880
$ws.Range("H130").Value = 28333.334
$ws.Range("J130").Value = 28333.334
$ws.Range("L130").Value = 28333.334
$ws.Range("N130").Value = -38373.334

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2188.4736
$ws.Range("I2").Value = 2313.1
$ws.Range("J2").Value = 2050
$ws.Range("K2").Value = 2313.1
$ws.Range("L2").Value = 2050
$ws.Range("M2").Value = -2200.1
$ws.Range("N2").Value = -2276
$ws.Range("H43").Value = 7188.5
$ws.Range("I43").Value = 4000
$ws.Range("K43").Value = 4000
$ws.Range("M43").Value = -3687
$ws.Range("H102").Value = 2004.4
$ws.Range("I102").Value = 1755.5
$ws.Range("K102").Value = 1755.5
$ws.Range("M102").Value = -133.5
$ws.Range("H116").Value = 2188.4736
$ws.Range("I116").Value = 2313.1
$ws.Range("J116").Value = 2050
$ws.Range("K116").Value = 2313.1
$ws.Range("L116").Value = 2050
$ws.Range("M116").Value = -19.09999999999991
$ws.Range("N116").Value = -6638
$ws.Range("H135").Value = 22982.166
$ws.Range("J135").Value = 22982.166
$ws.Range("L135").Value = 22982.166
$ws.Range("N135").Value = -33122.166

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2188.4736
$ws.Range("I3").Value = 2313.1
$ws.Range("J3").Value = 2050
$ws.Range("K3").Value = 2313.1
$ws.Range("L3").Value = 2050
$ws.Range("M3").Value = -2199.1
$ws.Range("N3").Value = -2278
$ws.Range("H86").Value = 2007.7693
$ws.Range("I86").Value = 1962.625
$ws.Range("J86").Value = 2080
$ws.Range("K86").Value = 1962.625
$ws.Range("L86").Value = 2080
$ws.Range("M86").Value = -839.625
$ws.Range("N86").Value = -4326
$ws.Range("H89").Value = 2007.7693
$ws.Range("I89").Value = 1962.625
$ws.Range("J89").Value = 2080
$ws.Range("K89").Value = 9813.125
$ws.Range("L89").Value = 10400
$ws.Range("M89").Value = -4197.125
$ws.Range("N89").Value = -21632
$ws.Range("H99").Value = 1921.2222
$ws.Range("I99").Value = 1622.5
$ws.Range("K99").Value = 1622.5
$ws.Range("M99").Value = -124.5
$ws.Range("H107").Value = 1831.5834
$ws.Range("I107").Value = 1694.125
$ws.Range("K107").Value = 1694.125
$ws.Range("M107").Value = 225.875
$ws.Range("H110").Value = 28500
$ws.Range("J110").Value = 28500
$ws.Range("L110").Value = 28500
$ws.Range("N110").Value = -36680

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2807.2666
$ws.Range("I132").Value = 1780.1666
$ws.Range("J132").Value = 4347.9165
$ws.Range("K132").Value = 5340.4998
$ws.Range("L132").Value = 13043.7495
$ws.Range("M132").Value = -2810.4998
$ws.Range("N132").Value = -18103.7495

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 882.3333
$ws.Range("I5").Value = 497.125
$ws.Range("J5").Value = 1190.5
$ws.Range("K5").Value = 1491.375
$ws.Range("L5").Value = 3571.5
$ws.Range("M5").Value = -1379.375
$ws.Range("N5").Value = -3795.5
$ws.Range("H94").Value = 3607.9583
$ws.Range("I94").Value = 1966.6666
$ws.Range("J94").Value = 3842.4285
$ws.Range("K94").Value = 5899.9998
$ws.Range("L94").Value = 11527.2855
$ws.Range("M94").Value = -5223.9998
$ws.Range("N94").Value = -12879.2855
$ws.Range("H135").Value = 882.3333
$ws.Range("I135").Value = 497.125
$ws.Range("J135").Value = 1190.5
$ws.Range("K135").Value = 4474.125
$ws.Range("L135").Value = 10714.5
$ws.Range("M135").Value = -1939.125
$ws.Range("N135").Value = -15784.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 68.25
$ws.Range("I2").Value = 41.916668
$ws.Range("J2").Value = 94.583336
$ws.Range("K2").Value = 41.916668
$ws.Range("L2").Value = 94.583336
$ws.Range("M2").Value = 71.083332
$ws.Range("N2").Value = -320.583336

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1623.3334
$ws.Range("J82").Value = 1185
$ws.Range("L82").Value = 1185
$ws.Range("N82").Value = -1907
$ws.Range("H85").Value = 1623.3334
$ws.Range("J85").Value = 1185
$ws.Range("L85").Value = 1185
$ws.Range("N85").Value = -3681
$ws.Range("H136").Value = 12346982
$ws.Range("I136").Value = 14493930
$ws.Range("J136").Value = 2025
$ws.Range("K136").Value = 43481790
$ws.Range("L136").Value = 6075
$ws.Range("M136").Value = -43479240
$ws.Range("N136").Value = -11175
$ws.Range("H141").Value = 39684
$ws.Range("J141").Value = 39684
$ws.Range("L141").Value = 39684
$ws.Range("N141").Value = -50044

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 30024.75
$ws.Range("J92").Value = 30024.75
$ws.Range("L92").Value = 30024.75
$ws.Range("N92").Value = -35016.75
$ws.Range("H136").Value = 20410220
$ws.Range("I136").Value = 35716836
$ws.Range("J136").Value = 1401.8096
$ws.Range("K136").Value = 107150508
$ws.Range("L136").Value = 4205.4288
$ws.Range("M136").Value = -107147958
$ws.Range("N136").Value = -9305.4288
